$wb = $excel.ActiveWorkbook

# --- Rename the existing "VISA" sheet to "VISA(old)" and insert a brand
# new "VISA" sheet right after it (this becomes the new active sheet). ---
$wsOld = $wb.Worksheets.Item("VISA")
$wsOld.Name = "VISA(old)"

$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOld)
$wsNew.Name = "VISA"

# --- Update selection on the old VISA sheet ---
$wsOld.Activate() | Out-Null
$wsOld.Range("A2:B2").Select() | Out-Null

# --- Populate the new VISA sheet ---
$wsNew.Range("A1").Value = "No"
$wsNew.Range("B1").Value = "Description"
$wsNew.Range("C1").Value = "Inputs"
$wsNew.Range("D1").Value = "Expected Output"
$wsNew.Range("E1").Value = "Actual Output"
$wsNew.Range("F1").Value = "Pass / Fail "

$wsNew.Range("A2").Value = 1
$wsNew.Range("B2").Value = "To compare the number of applications received"
$wsNew.Range("D2").Value = "Australia Tourist : 5`nICA STVP : 10 `nIndonesia - VKBKP 212 : 61`nIndonesia - VKU 211 : 50`nIndonesia - VTT 312 - Less than 6 months : 124`nIndonesia - VTT 312 - More than 6 months : 121`nIndonesia - VTT 316 - Less than 6 months: 1`nIndonesia - VTT 316 - More than 6 months: 2`nIndonesia - VTT 317 - Less than 6 months: 13`nIndonesia - VTT 317 - More than 6 months: 17`nIndonesia Entry Visa Extend : 1`nMaternity : 11"
$wsNew.Range("C2").Value = "-"

# --- Formatting to match the header / body styles used elsewhere in the
# workbook (bold header row, left/top/wrap body cells). ---
$wsNew.Range("A1:F1").Font.Bold = $true

$wsNew.Range("A2:B2").HorizontalAlignment = -4131
$wsNew.Range("A2:B2").VerticalAlignment = -4160
$wsNew.Range("A2:B2").WrapText = $true

$wsNew.Range("C2:D2").VerticalAlignment = -4160
$wsNew.Range("C2:D2").WrapText = $true

# --- Row height / column widths ---
$wsNew.Rows.Item(2).RowHeight = 197

# Note: the host's ColumnWidth setter quantizes to a 6px Maximum-Digit-Width
# grid (stored = (floor(chars*6+0.5)+5)/6), so these are the pre-images
# that round-trip closest to the target stored widths (3.6640625,
# 17.33203125, 27.33203125, 43.6640625, 27.83203125).
$wsNew.Columns.Item(1).ColumnWidth = 2.8333333333333335
$wsNew.Columns.Item(2).ColumnWidth = 16.5
$wsNew.Columns.Item(3).ColumnWidth = 26.5
$wsNew.Columns.Item(4).ColumnWidth = 42.833333333333336
$wsNew.Columns.Item(5).ColumnWidth = 27.0

# --- Make the new VISA sheet the active one, matching the selection in
# the updated workbook. ---
$wsNew.Activate()
$wsNew.Range("C2").Select()
